# Rename the third worksheet ("CreateAccountTest") to "SearchTest"
# and update the related test-suite row + selection/active-tab state,
# matching the "Files updated to run with Extent Reports" commit.

$wb = $excel.ActiveWorkbook

$suiteSheet = $wb.Worksheets.Item("Test_Suite")
$createSheet = $wb.Worksheets.Item("CreateAccountTest")

# Rename the sheet itself.
$createSheet.Name = "SearchTest"

# The Test_Suite row that referenced the old sheet name needs to point at
# the new one (the "y" value in column B is unchanged).
$suiteSheet.Range("A3").Value = "SearchTest"
$suiteSheet.Range("B3").Value = "y"

# Update the selection on the renamed sheet and make it the active tab.
$createSheet.Activate()
$createSheet.Range("H16").Select()
